# Lowercase the category labels in columns C, D, E, F for rows 2 and 3
# (gender, location, mood, activity), which were previously capitalized.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "male"
$ws.Range("D2").Value = "outdoor"
$ws.Range("E2").Value = "excited"
$ws.Range("F2").Value = "working out"

$ws.Range("C3").Value = "male"
$ws.Range("D3").Value = "home"
$ws.Range("E3").Value = "relaxed"
$ws.Range("F3").Value = "amusing"

# Update the active selection to match the saved workbook state (C4).
$ws.Range("C4").Select()
